$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Swap the content (columns B:V) of rows 22 and 23.
#    Column A (Indice) keeps its original per-row value.
# ---------------------------------------------------------------
$row22 = @{
    B = "belgium"; C = "jupiler-pro-league"; D = "2023-2024"
    F = "Eupen"; G = 0; H = "Club Brugge KV"; I = 5
    J = 3.76; K = "06/08/2023 18:42"
    L = 8.35; M = "13/08/2023 15:59"
    N = 4.19; O = "06/08/2023 18:42"
    P = 5.47; Q = "13/08/2023 15:59"
    R = 1.85; S = "06/08/2023 18:42"
    T = 1.36; U = "13/08/2023 15:59"
    V = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/eupen-club-brugge/lWaC8fiN/"
}

$row23 = @{
    B = "belgium"; C = "jupiler-pro-league"; D = "2023-2024"
    F = "Cercle Brugge KSV"; G = 0; H = "Genk"; I = 1
    J = 2.55; K = "06/08/2023 16:12"
    L = 2.52; M = "13/08/2023 15:56"
    N = 3.75; O = "06/08/2023 16:12"
    P = 3.74; Q = "13/08/2023 15:59"
    R = 2.49; S = "06/08/2023 16:12"
    T = 2.71; U = "13/08/2023 15:56"
    V = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/cercle-brugge-genk/SGB3AGMA/"
}

foreach ($col in $row22.Keys) {
    $ws.Range(($col + "22")).Value = $row22[$col]
}
foreach ($col in $row23.Keys) {
    $ws.Range(($col + "23")).Value = $row23[$col]
}

# ---------------------------------------------------------------
# 2) Swap the content (columns B:V) of rows 30 and 31.
# ---------------------------------------------------------------
$row30 = @{
    B = "belgium"; C = "jupiler-pro-league"; D = "2023-2024"
    F = "Gent"; G = 2; H = "St. Truiden"; I = 2
    J = 1.4; K = "13/08/2023 19:42"
    L = 1.52; M = "20/08/2023 15:57"
    N = 5; O = "13/08/2023 19:42"
    P = 4.45; Q = "20/08/2023 15:59"
    R = 6.31; S = "13/08/2023 19:42"
    T = 6.41; U = "20/08/2023 15:59"
    V = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/gent-st-truiden/SxmvNg6G/"
}

$row31 = @{
    B = "belgium"; C = "jupiler-pro-league"; D = "2023-2024"
    F = "Genk"; G = 0; H = "Charleroi"; I = 0
    J = 1.47; K = "13/08/2023 18:42"
    L = 1.59; M = "20/08/2023 15:53"
    N = 4.85; O = "13/08/2023 18:42"
    P = 4.57; Q = "20/08/2023 15:53"
    R = 5.52; S = "13/08/2023 18:42"
    T = 5.25; U = "20/08/2023 15:58"
    V = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/genk-charleroi/xlfRPeMc/"
}

foreach ($col in $row30.Keys) {
    $ws.Range(($col + "30")).Value = $row30[$col]
}
foreach ($col in $row31.Keys) {
    $ws.Range(($col + "31")).Value = $row31[$col]
}

# ---------------------------------------------------------------
# 3) Append two new rows (146 and 147) after the previous last
#    row (145), extending the sheet dimension to A1:V147.
#    Copy the formatting (styles) from row 145 first, then fill
#    in the values.
# ---------------------------------------------------------------
$ws.Range("A145:V145").Copy()
$ws.Range("A146:V146").PasteSpecial(-4122)
$ws.Range("A147:V147").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$row146 = @{
    A = 145; B = "belgium"; C = "jupiler-pro-league"; D = "2023-2024"
    E = 45280.77083333334
    F = "Cercle Brugge KSV"; G = 3; H = "Kortrijk"; I = 0
    J = 1.36; K = "17/12/2023 19:43"
    L = 1.29; M = "20/12/2023 18:29"
    N = 5.11; O = "17/12/2023 19:43"
    P = 5.91; Q = "20/12/2023 18:29"
    R = 7.32; S = "17/12/2023 19:43"
    T = 10.37; U = "20/12/2023 18:29"
    V = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/cercle-brugge-kortrijk/823Cx9ui/"
}

$row147 = @{
    A = 146; B = "belgium"; C = "jupiler-pro-league"; D = "2023-2024"
    E = 45280.86458333334
    F = "KV Mechelen"; G = 3; H = "St. Liege"; I = 0
    J = 2.54; K = "17/12/2023 16:12"
    L = 2.64; M = "20/12/2023 20:42"
    N = 3.33; O = "17/12/2023 16:12"
    P = 3.14; Q = "20/12/2023 20:18"
    R = 2.7; S = "17/12/2023 16:12"
    T = 2.97; U = "20/12/2023 20:42"
    V = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/kv-mechelen-st-liege/ObgyX49S/"
}

foreach ($col in $row146.Keys) {
    $ws.Range(($col + "146")).Value = $row146[$col]
}
foreach ($col in $row147.Keys) {
    $ws.Range(($col + "147")).Value = $row147[$col]
}
